# Atualização automática: 2025-08-11 08:49:38
# Cyclically shift the data in rows 7-11 (columns A, D, E, F, G, H, I, J) up
# by one row: row 8's values move to row 7, row 9's to row 8, row 10's to
# row 9, row 11's to row 10, and the original row 7's values wrap to row 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "D", "E", "F", "G", "H", "I", "J")
# I and J hold numeric-looking text ("962,713,1006,765", "0.76") that Excel
# would otherwise auto-convert to a number (stripping commas / changing
# type). Force those columns to Text format so the written values keep
# their original string representation.
$textCols = @("I", "J")

# Capture the original values for rows 7..11 before any writes.
$orig = @{}
foreach ($r in 7..11) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value2
    }
}

foreach ($c in $textCols) {
    $ws.Range("$c" + "7:" + "$c" + "11").NumberFormat = "@"
}

# Row 7 <- Row 8, Row 8 <- Row 9, Row 9 <- Row 10, Row 10 <- Row 11
foreach ($r in 7..10) {
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $orig[$r + 1][$c]
    }
}

# Row 11 <- original Row 7
foreach ($c in $cols) {
    $ws.Range("$c" + "11").Value2 = $orig[7][$c]
}
